# Update odds/values on Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.62

# Row 4
$ws.Range("H4").Value = 3.6
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 6
$ws.Range("N4").Value = 7.5
$ws.Range("AJ4").Value = 19
$ws.Range("AZ4").Value = 126
$ws.Range("BA4").Value = 151

# Row 5
$ws.Range("G5").Value = 1.47
$ws.Range("H5").Value = 3.55
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 1.93
$ws.Range("K5").Value = 2.18
$ws.Range("L5").Value = 7.2
$ws.Range("N5").Value = 6.65
$ws.Range("O5").Value = 1.34
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.65
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.57
$ws.Range("W5").Value = 5.2
$ws.Range("X5").Value = 5.9
$ws.Range("Y5").Value = 8.25
$ws.Range("Z5").Value = 9.75
$ws.Range("AA5").Value = 13.5
$ws.Range("AC5").Value = 7.9
$ws.Range("AD5").Value = 7.3
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 120
$ws.Range("AH5").Value = 17
$ws.Range("AI5").Value = 55
$ws.Range("AJ5").Value = 25
$ws.Range("AK5").Value = 250
$ws.Range("AL5").Value = 110
$ws.Range("AM5").Value = 100
$ws.Range("AN5").Value = 3.15
$ws.Range("AO5").Value = 6.6
$ws.Range("AP5").Value = 16.5
$ws.Range("AR5").Value = 50
$ws.Range("AS5").Value = 250
$ws.Range("AT5").Value = 2.57
$ws.Range("AU5").Value = 7.9
$ws.Range("AV5").Value = 75
$ws.Range("AW5").Value = 8.75
$ws.Range("AX5").Value = 50
$ws.Range("AY5").Value = 45
$ws.Range("AZ5").Value = 400
$ws.Range("BA5").Value = 350

$wb.Save()
